# Refresh coin Price (D) and Volume(1h) (E) figures to match
# the latest scrape (GitHub Actions cron refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.550.86"
$ws.Range("E2").Value = "  +0.34%  "
$ws.Range("D3").Value = "1.816.64"
$ws.Range("E3").Value = "  +0.51%  "
$ws.Range("E4").Value = "  -0.58%  "
$ws.Range("D5").Formula = "'1.002"
$ws.Range("E5").Value = "  -0.57%  "
$ws.Range("D6").Formula = "'306.25"
$ws.Range("E6").Value = "  -0.65%  "
$ws.Range("D7").Formula = "'0.4532"
$ws.Range("E7").Value = "  -0.67%  "
$ws.Range("D8").Formula = "'0.3596"
$ws.Range("E8").Value = "  -1.81%  "
$ws.Range("D9").Formula = "'46.33"
$ws.Range("E9").Value = "  +2.89%  "
$ws.Range("D10").Formula = "'0.07091"
$ws.Range("E10").Value = "  -0.63%  "
$ws.Range("D11").Formula = "'0.8963"
$ws.Range("E11").Value = "  +2.21%  "
$ws.Range("D12").Formula = "'0.07768"
$ws.Range("E12").Value = "  -0.10%  "
$ws.Range("D13").Formula = "'19.36"
$ws.Range("E13").Value = "  +0.17%  "
$ws.Range("D14").Value = "1.834.96"
$ws.Range("E14").Value = "  +2.46%  "
$ws.Range("D15").Formula = "'5.268"
$ws.Range("E15").Value = "  -0.17%  "
$ws.Range("D16").Formula = "'6.322"
$ws.Range("E16").Value = "  -0.81%  "
$ws.Range("D17").Formula = "'85.33"
$ws.Range("E17").Value = "  -0.66%  "
$ws.Range("D18").Formula = "'1.004"
$ws.Range("E18").Value = "  -0.54%  "
$ws.Range("D19").Formula = "'0.000008602"
$ws.Range("E19").Value = "  +0.27%  "
$ws.Range("E20").Value = "  -0.51%  "
$ws.Range("D21").Value = "26.591.80"
$ws.Range("E21").Value = "  +0.35%  "
$ws.Range("E22").Value = "  -0.52%  "
$ws.Range("D23").Formula = "'4.964"
$ws.Range("E23").Value = "  -0.56%  "
$ws.Range("D24").Formula = "'10.52"
$ws.Range("E24").Value = "  +0.66%  "
$ws.Range("D25").Formula = "'1.968"
$ws.Range("E25").Value = "  -0.72%  "
$ws.Range("D26").Formula = "'151.17"
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("E27").Value = "  -0.87%  "
$ws.Range("D28").Formula = "'2.050"
$ws.Range("E28").Value = "  +0.26%  "
$ws.Range("D29").Formula = "'112.55"
$ws.Range("E30").Value = "  +0.24%  "
$ws.Range("D31").Formula = "'0.08726"
$ws.Range("E31").Value = "  +0.63%  "
$ws.Range("D32").Formula = "'3.123"
$ws.Range("E32").Value = "  +2.69%  "
$ws.Range("D33").Formula = "'0.7580"
$ws.Range("E33").Value = "  +3.86%  "
$ws.Range("D34").Formula = "'2.734"
$ws.Range("E34").Value = "  +9.84%  "
$ws.Range("D35").Formula = "'4.433"
$ws.Range("E35").Value = "  -0.53%  "
$ws.Range("D36").Formula = "'1.115"
$ws.Range("E36").Value = "  +0.14%  "
$ws.Range("D37").Formula = "'1.073"
$ws.Range("E37").Value = "  -0.52%  "
$ws.Range("D38").Formula = "'0.01938"
$ws.Range("E38").Value = "  +0.30%  "
$ws.Range("D39").Formula = "'2.906"
$ws.Range("E39").Value = "  +0.54%  "
$ws.Range("D40").Formula = "'0.05106"
$ws.Range("E40").Value = "  +0.04%  "
$ws.Range("D41").Formula = "'0.5110"
$ws.Range("E41").Value = "  +1.95%  "
$ws.Range("D42").Formula = "'6.774"
$ws.Range("E42").Value = "  -2.67%  "
$ws.Range("D43").Formula = "'0.1512"
$ws.Range("E43").Value = "  -2.91%  "
$ws.Range("D44").Formula = "'8.046"
$ws.Range("E44").Value = "  -0.97%  "
$ws.Range("D45").Formula = "'0.4706"
$ws.Range("E45").Value = "  +2.19%  "
$ws.Range("D46").Formula = "'1.002"
$ws.Range("E46").Value = "  -0.66%  "
$ws.Range("D47").Formula = "'10.02"
$ws.Range("E47").Value = "  +0.61%  "
$ws.Range("D48").Formula = "'101.17"
$ws.Range("E48").Value = "  +0.19%  "
$ws.Range("D49").Formula = "'1.577"
$ws.Range("E49").Value = "  -0.73%  "
$ws.Range("D50").Formula = "'0.05984"
$ws.Range("E50").Value = "  -0.51%  "
$ws.Range("D51").Formula = "'63.92"
$ws.Range("E51").Value = "  -0.22%  "
